$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark exercises 1, 2 and 3 as done by setting these cells to 1
$ws.Range("O2:Q4").Value = 1

# Update the view: scroll so row 22 is at the top, and select Q4
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q4").Select()
